# Add a "backend" attribute column (with value "varken") to the
# "entities" sheet, and leave the "entities" sheet as the active /
# selected sheet with C3 selected (mirrors the author's edit in Excel).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("entities")

$ws.Range("C1").Value = "backend"
$ws.Range("C2").Value = "varken"

$ws.Activate()
$ws.Range("C3").Select()
